# Update crypto price/volume table per GitHub Actions scrape run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value2 = "26.899.35"
$ws.Cells.Item(2, 5).Value2 = "  +0.23%  "

$ws.Cells.Item(3, 4).Value2 = "1.548.09"
$ws.Cells.Item(3, 5).Value2 = "  -0.87%  "

$ws.Cells.Item(4, 5).Value2 = "  +0.26%  "

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value2 = "206.98"
$ws.Cells.Item(5, 5).Value2 = "  +0.60%  "

$ws.Cells.Item(6, 5).Value2 = "  +0.02%  "

$ws.Cells.Item(7, 5).Value2 = "  +0.25%  "

$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value2 = "0.247"
$ws.Cells.Item(8, 5).Value2 = "  -0.06%  "

$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value2 = "21.44"
$ws.Cells.Item(9, 5).Value2 = "  -1.20%  "

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value2 = "0.0583"
$ws.Cells.Item(10, 5).Value2 = "  +0.00%  "

$ws.Cells.Item(11, 5).Value2 = "  -0.98%  "

$ws.Cells.Item(12, 4).Value2 = "1.768.53"
$ws.Cells.Item(12, 5).Value2 = "  -0.84%  "

$ws.Cells.Item(13, 4).Value2 = "1.549.94"
$ws.Cells.Item(13, 5).Value2 = "  -0.62%  "

$ws.Cells.Item(14, 5).Value2 = "  -0.56%  "

$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value2 = "0.512"
$ws.Cells.Item(15, 5).Value2 = "  -0.05%  "

$ws.Cells.Item(16, 4).Value2 = "26.906.96"
$ws.Cells.Item(16, 5).Value2 = "  +0.26%  "

$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value2 = "61.44"
$ws.Cells.Item(17, 5).Value2 = "  +0.59%  "

$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value2 = "214.64"
$ws.Cells.Item(18, 5).Value2 = "  +0.34%  "

$ws.Cells.Item(19, 5).Value2 = "  +1.03%  "

$ws.Cells.Item(20, 5).Value2 = "  -1.54%  "

$ws.Cells.Item(21, 5).Value2 = "  +0.20%  "

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value2 = "4.02"
$ws.Cells.Item(22, 5).Value2 = "  -2.31%  "

$ws.Cells.Item(23, 5).Value2 = "  +0.03%  "

$ws.Cells.Item(24, 5).Value2 = "  -2.70%  "

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value2 = "151.93"
$ws.Cells.Item(25, 5).Value2 = "  -1.21%  "

$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value2 = "6.63"
$ws.Cells.Item(26, 5).Value2 = "  -1.10%  "

$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value2 = "14.88"
$ws.Cells.Item(27, 5).Value2 = "  -0.30%  "

$ws.Cells.Item(28, 5).Value2 = "  +0.26%  "

$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value2 = "0.104"
$ws.Cells.Item(29, 5).Value2 = "  +0.66%  "

$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value2 = "0.0459"
$ws.Cells.Item(30, 5).Value2 = "  -0.56%  "

$ws.Cells.Item(31, 5).Value2 = "  -0.36%  "

$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value2 = "3.23"
$ws.Cells.Item(32, 5).Value2 = "  +2.36%  "

$ws.Cells.Item(33, 4).Value2 = "1.367.97"
$ws.Cells.Item(33, 5).Value2 = "  -2.29%  "

$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value2 = "2.96"
$ws.Cells.Item(34, 5).Value2 = "  +1.60%  "

$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value2 = "1.53"
$ws.Cells.Item(35, 5).Value2 = "  +0.57%  "

$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value2 = "0.959"
$ws.Cells.Item(36, 5).Value2 = "  +4.24%  "

$ws.Cells.Item(37, 5).Value2 = "  +0.36%  "

$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value2 = "0.0165"
$ws.Cells.Item(38, 5).Value2 = "  -0.16%  "

$ws.Cells.Item(39, 5).Value2 = "  -0.47%  "

$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value2 = "0.806"
$ws.Cells.Item(40, 5).Value2 = "  -0.67%  "

$ws.Cells.Item(41, 5).Value2 = "  +0.25%  "

$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value2 = "5.63"
$ws.Cells.Item(42, 5).Value2 = "  +6.21%  "

$ws.Cells.Item(43, 5).Value2 = "  -0.16%  "

$ws.Cells.Item(44, 5).Value2 = "  +2.18%  "

$ws.Cells.Item(45, 5).Value2 = "  +1.09%  "

$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value2 = "1.73"
$ws.Cells.Item(46, 5).Value2 = "  -1.71%  "

$ws.Cells.Item(47, 2).Value2 = "RocketPoolETH"
$ws.Cells.Item(47, 3).Value2 = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Cells.Item(47, 4).Value2 = "1.683.24"
$ws.Cells.Item(47, 5).Value2 = "  -0.75%  "

$ws.Cells.Item(48, 2).Value2 = "Quant"
$ws.Cells.Item(48, 3).Value2 = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value2 = "85.61"
$ws.Cells.Item(48, 5).Value2 = "  -0.42%  "

$ws.Cells.Item(49, 2).Value2 = "Cronos"
$ws.Cells.Item(49, 3).Value2 = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value2 = "0.0508"
$ws.Cells.Item(49, 5).Value2 = "  +1.03%  "

$ws.Cells.Item(50, 2).Value2 = "BabyDogeCoin"
$ws.Cells.Item(50, 3).Value2 = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Cells.Item(50, 4).Value2 = "0.0₇0974"
$ws.Cells.Item(50, 5).Value2 = "  -1.08%  "

$ws.Cells.Item(51, 2).Value2 = "Algorand"
$ws.Cells.Item(51, 3).Value2 = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value2 = "0.0949"
$ws.Cells.Item(51, 5).Value2 = "  +0.56%  "
